# Add a new date column "04-ago" to the right of the existing last
# column ("03-ago"), shifting that day's values into the new column
# and filling the former last column with the new day's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column AS (column 45), continuing the "03-ago" series.
$ws.Cells.Item(1, 45).Value = "04-ago"

# New values for the new "04-ago" column (AS) -- these are the figures
# that used to live in the "03-ago" column (AR) before the new day was
# inserted.
$asValues = @{
    2  = 0
    3  = 15.689327096616635
    4  = 19.966936274338202
    5  = 18.822531355581052
    6  = 0
    7  = 15.114920988675921
    8  = 8.3982846636579307
    9  = 12.541465128871852
    10 = 11.648992059434557
    11 = 15.023644843967348
    12 = 0
    13 = 6.6675980200815719
    14 = 0
    15 = 0
    16 = 11.506558368047143
    17 = 0
    18 = 0
}

# Updated values for the existing "03-ago" column (AR) now that the new
# day's data has arrived.
$arValues = @{
    3  = 16.90502550478282
    4  = 18.25753844709287
    5  = 20.080904446517032
    7  = 13.476089988166191
    8  = 7.6357437427157979
    9  = 13.755462932008793
    10 = 14.56767543451641
    11 = 14.188463642094328
    13 = 8.7009517800518097
    16 = 11.778161198341493
}

foreach ($row in 2..18) {
    $ws.Cells.Item($row, 45).Value = $asValues[$row]
    if ($arValues.ContainsKey($row)) {
        $ws.Cells.Item($row, 44).Value = $arValues[$row]
    }
}

# Move the active selection one column to the right, matching the
# author's cursor position after inserting the new column.
$ws.Range("AU7").Select()
